# Append bank-account hints to the five bank-journal names in column C.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "B. Pop. Software (IT15*456)"
$ws.Range("C3").Value = "B. Credito per Tutti (IT74*680)"
$ws.Range("C4").Value = "Portafoglio RiBA B. Pop. Soft. (IT15*456/IT26*456)"
$ws.Range("C5").Value = "Portafoglio Anticipi B. Pop. Soft. (IT15*456/IT82*456)"
$ws.Range("C6").Value = "Portafoglio Anticipi Estero BCT (IT74*680/IT58*680)"

# Widen column C so the longer names fit (was ~32.55 characters, now ~46 characters).
$ws.Columns.Item(3).ColumnWidth = 45.33

# Move the active selection to C6, matching the cell last edited.
$ws.Range("C6").Select()
